$d = $word.ActiveDocument

$d.Content.Find.Execute("69÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=", 2) | Out-Null
$d.Content.Find.Execute("82÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷4=", 2) | Out-Null
$d.Content.Find.Execute("75÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷5=", 2) | Out-Null
$d.Content.Find.Execute("60÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷8=", 2) | Out-Null
$d.Content.Find.Execute("82÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=", 2) | Out-Null
$d.Content.Find.Execute("17÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷9=", 2) | Out-Null
$d.Content.Find.Execute("82÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷6=", 2) | Out-Null
$d.Content.Find.Execute("46÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷4=", 2) | Out-Null
$d.Content.Find.Execute("58÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷8=", 2) | Out-Null
$d.Content.Find.Execute("14÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷3=", 2) | Out-Null
$d.Content.Find.Execute("90÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷8=", 2) | Out-Null
$d.Content.Find.Execute("21÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷2=", 2) | Out-Null
$d.Content.Find.Execute("99÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷9=", 2) | Out-Null
$d.Content.Find.Execute("56÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=", 2) | Out-Null
$d.Content.Find.Execute("83÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷8=", 2) | Out-Null
$d.Content.Find.Execute("87÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=", 2) | Out-Null
$d.Content.Find.Execute("44÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷7=", 2) | Out-Null
$d.Content.Find.Execute("31÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=", 2) | Out-Null
$d.Content.Find.Execute("89÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=", 2) | Out-Null
$d.Content.Find.Execute("22÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷4=", 2) | Out-Null
$d.Content.Find.Execute("97÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷8=", 2) | Out-Null
$d.Content.Find.Execute("52÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷2=", 2) | Out-Null
$d.Content.Find.Execute("13÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=", 2) | Out-Null
$d.Content.Find.Execute("71÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷7=", 2) | Out-Null
$d.Content.Find.Execute("88÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷9=", 2) | Out-Null
